$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "2024-10-09T12:00+00:00", 3811062, 9, 2352.8, 3),
    @(3, "2024-10-09T13:00+00:00", 4694910, 9.199999999999999, 2352.8, 3.3),
    @(4, "2024-10-09T14:00+00:00", 5489148, 9.300000000000001, 2352.8, 3.2),
    @(5, "2024-10-09T15:00+00:00", 6156986, 9, 2352.8, 3.1),
    @(6, "2024-10-09T16:00+00:00", 6441854, 8.4, 2352.8, 3.2),
    @(7, "2024-10-09T17:00+00:00", 6456946, 8.300000000000001, 2352.8, 3.3),
    @(8, "2024-10-09T18:00+00:00", 6456946, 8.1, 2352.8, 2.6),
    @(9, "2024-10-09T19:00+00:00", 6456946, 8.6, 2352.8, 3.1),
    @(10, "2024-10-09T20:00+00:00", 6456946, 8.9, 2352.8, 3.3),
    @(11, "2024-10-09T21:00+00:00", 6456946, 9.1, 2352.8, 3),
    @(12, "2024-10-09T22:00+00:00", 6456946, 9.199999999999999, 2352.8, 4.4),
    @(13, "2024-10-09T23:00+00:00", 6456946, 9.5, 2352.8, 4.9),
    @(14, "2024-10-10T00:00+00:00", 6456946, 9.1, 2352.8, 3.7),
    @(15, "2024-10-10T01:00+00:00", 6456946, 9.1, 2352.8, 3.9),
    @(16, "2024-10-10T02:00+00:00", 6456946, 9, 2352.8, 4.2),
    @(17, "2024-10-10T03:00+00:00", 6456946, 8.6, 2352.8, 3.2),
    @(18, "2024-10-10T04:00+00:00", 6456946, 8.5, 2352.8, 2.6),
    @(19, "2024-10-10T05:00+00:00", 6456946, 8.4, 2352.8, 2.8),
    @(20, "2024-10-10T06:00+00:00", 6480528, 8.4, 2352.8, 2),
    @(21, "2024-10-10T07:00+00:00", 6608814, 8.9, 2352.8, 2.6),
    @(22, "2024-10-10T08:00+00:00", 7062528, 9.199999999999999, 2352.8, 1.3),
    @(23, "2024-10-10T09:00+00:00", 7765268, 9.300000000000001, 2352.8, 1.8),
    @(24, "2024-10-10T10:00+00:00", 8748160, 9.6, 2352.8, 2),
    @(25, "2024-10-10T11:00+00:00", 9038689.6, 9.800000000000001, 2352.8, 1.9),
    @(26, "2024-10-10T12:00+00:00", 9345253.6, 10.8, 2352.8, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
